$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that need to remain literal text (would otherwise be auto-converted
# to numbers/percentages by Excel) get NumberFormat "@" applied first.
$textCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E5",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E24",
    "D25",
    "E25",
    "D26",
    "E26",
    "E27",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "E47",
    "D48",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell.
$ws.Range('D2').Value = '245.38'
$ws.Range('E2').Value = '-0.56%'
$ws.Range('D3').Value = '28.29'
$ws.Range('E3').Value = '-4.83%'
$ws.Range('D4').Value = '5.249'
$ws.Range('E5').Value = '-0.50%'
$ws.Range('E6').Value = '0.70%'
$ws.Range('D7').Value = '3.201'
$ws.Range('E7').Value = '3.34%'
$ws.Range('D8').Value = '0.8508'
$ws.Range('E8').Value = '-0.72%'
$ws.Range('D9').Value = '0.9069'
$ws.Range('E9').Value = '4.79%'
$ws.Range('D10').Value = '0.1370'
$ws.Range('E10').Value = '0.35%'
$ws.Range('E11').Value = '0.03%'
$ws.Range('D12').Value = '0.03185'
$ws.Range('E12').Value = '8.79%'
$ws.Range('D13').Value = '0.09221'
$ws.Range('E13').Value = '-1.72%'
$ws.Range('D14').Value = '0.001528'
$ws.Range('E14').Value = '0.97%'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').Value = '0.005915'
$ws.Range('E15').Value = '-1.72%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').Value = '3.490'
$ws.Range('E16').Value = '0.00%'
$ws.Range('B17').Value = 'BTSEToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D17').Value = '2.175'
$ws.Range('E17').Value = '-4.28%'
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D18').Value = '0.01003'
$ws.Range('E18').Value = '1,577.69%'
$ws.Range('D19').Value = '0.3170'
$ws.Range('E19').Value = '-0.39%'
$ws.Range('D20').Value = '0.03307'
$ws.Range('E20').Value = '-2.54%'
$ws.Range('D21').Value = '0.1276'
$ws.Range('E21').Value = '-2.04%'
$ws.Range('D22').Value = '3.523'
$ws.Range('E22').Value = '1.54%'
$ws.Range('D23').Value = '0.04074'
$ws.Range('E24').Value = '-0.09%'
$ws.Range('D25').Value = '0.001220'
$ws.Range('E25').Value = '-0.32%'
$ws.Range('D26').Value = '0.004157'
$ws.Range('E26').Value = '-16.99%'
$ws.Range('E27').Value = '-0.87%'
$ws.Range('D40').Value = '0.03793'
$ws.Range('E40').Value = '1.10%'
$ws.Range('D41').Value = '0.1065'
$ws.Range('E41').Value = '-0.59%'
$ws.Range('D42').Value = '0.003735'
$ws.Range('E42').Value = '-35.24%'
$ws.Range('D43').Value = '0.002198'
$ws.Range('E43').Value = '-7.75%'
$ws.Range('E44').Value = '7.55%'
$ws.Range('D45').Value = '0.00005260'
$ws.Range('E45').Value = '0.15%'
$ws.Range('D46').Value = '0.00000000749'
$ws.Range('E46').Value = '-0.04%'
$ws.Range('E47').Value = '62.22%'
$ws.Range('D48').Value = '0.002268'
$ws.Range('E48').Value = '0.40%'
$ws.Range('D49').Value = '0.00002098'
$ws.Range('E49').Value = '-0.04%'
$ws.Range('D50').Value = '0.0001998'
$ws.Range('E50').Value = '-0.04%'
